# Update of results and scripts. Anonymized "fedcore" -> "approach", and
# give the C1/D1 (and F1/G1 on sheet2) merged-header top cells their own
# border treatment (top+bottom, and top+right+bottom for the rightmost).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Rename "fedcore" header label to "approach" wherever it appears in row 2.
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value -eq "fedcore") {
            $cell.Value = "approach"
        }
    }
}

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws1.Range("C1").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws1.Range("C1").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws1.Range("D1").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws1.Range("D1").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws1.Range("D1").Borders.Item(10).LineStyle = 1  # xlEdgeRight

$ws2 = $wb.Worksheets.Item("computational_comparison")
foreach ($addr in @("C1", "F1")) {
    $ws2.Range($addr).Borders.Item(8).LineStyle = 1
    $ws2.Range($addr).Borders.Item(9).LineStyle = 1
}
foreach ($addr in @("D1", "G1")) {
    $ws2.Range($addr).Borders.Item(8).LineStyle = 1
    $ws2.Range($addr).Borders.Item(9).LineStyle = 1
    $ws2.Range($addr).Borders.Item(10).LineStyle = 1
}

# G5 on computational_comparison was an empty inline-string cell; clear it
# fully so it no longer exists as a distinct cell entry.
$ws2.Range("G5").ClearContents()
